# Fix typo in the project-title header cell and update the selection/view state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix accent typo: "Titúlo del proyecto" -> "Título del proyecto"
$ws.Range("C2").Value = "Título del proyecto"

# Update sheet view: remove custom topLeftCell, move selection to C6
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C6").Select()
